# Swap the order of "Recorded By" names in column G from
# "dnasr281@gmail.com, System" to "System, dnasr281@gmail.com"
# for every row in the used range of the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$used = $ws.UsedRange
$lastRow = $used.Rows().Count()

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()
    if ($val -eq $oldValue) {
        $cell.Value = $newValue
    }
}
